$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: "static function" row.
# C8 changes from "x" to the checkmark used elsewhere in the sheet, and
# picks up the same border/alignment formatting that D8 already has.
$d8 = $ws.Range("D8")
$c8 = $ws.Range("C8")
$d8.Copy()
$c8.PasteSpecial(-4122)
$c8.Value = [char]0x2713

# F8's note ("Field::Type gets set to decltype(&T::field)") is cleared out.
$ws.Range("F8").Value = ""

# A new blank, centered, Arial-styled cell appears at F10 (formatting-only row).
# Borrow the centered Arial heading format already used elsewhere (e.g. C3)
# and strip its border so F10 ends up border-free.
$c3 = $ws.Range("C3")
$f10 = $ws.Range("F10")
$c3.Copy()
$f10.PasteSpecial(-4122)
$f10.Borders.LineStyle = -4142

# Selection moves to F14.
$ws.Range("F14").Select()
